$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H7").Value = 3311.111
$ws.Range("I7").Value = 0
$ws.Range("J7").Value = 3311.111
$ws.Range("K7").Value = 0
$ws.Range("L7").Value = 3311.111
$ws.Range("M7").ClearContents()
$ws.Range("N7").Value = -3535.111
$ws.Range("H14").Value = 3311.111
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 3311.111
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 3311.111
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -3693.111
$ws.Range("H53").Value = 338.8
$ws.Range("I53").Value = 330.64285
$ws.Range("J53").Value = 357.83334
$ws.Range("K53").Value = 330.64285
$ws.Range("L53").Value = 357.83334
$ws.Range("M53").Value = 306.35715
$ws.Range("N53").Value = -1631.83334
$ws.Range("H64").Value = 7998
$ws.Range("I64").Value = 11330
$ws.Range("J64").Value = 3000
$ws.Range("K64").Value = 11330
$ws.Range("L64").Value = 3000
$ws.Range("M64").Value = -11082
$ws.Range("N64").Value = -3496
$ws.Range("H67").Value = 7998
$ws.Range("I67").Value = 11330
$ws.Range("J67").Value = 3000
$ws.Range("K67").Value = 11330
$ws.Range("L67").Value = 3000
$ws.Range("M67").Value = -10472
$ws.Range("N67").Value = -4716
$ws.Range("H74").Value = 6493140.5
$ws.Range("I74").Value = 7420160.5
$ws.Range("J74").Value = 4000
$ws.Range("K74").Value = 7420160.5
$ws.Range("L74").Value = 4000
$ws.Range("M74").Value = -7419224.5
$ws.Range("N74").Value = -5872
$ws.Range("H77").Value = 6493140.5
$ws.Range("I77").Value = 7420160.5
$ws.Range("J77").Value = 4000
$ws.Range("K77").Value = 37100802.5
$ws.Range("L77").Value = 20000
$ws.Range("M77").Value = -37096122.5
$ws.Range("N77").Value = -29360
$ws.Range("H100").Value = 1000
$ws.Range("I100").Value = 1000
$ws.Range("J100").Value = 0
$ws.Range("K100").Value = 1000
$ws.Range("L100").Value = 0
$ws.Range("M100").Value = -459
$ws.Range("N100").ClearContents()
$ws.Range("H129").Value = 841820.9399999999
$ws.Range("I129").Value = 255.36363
$ws.Range("J129").Value = 1092016.1
$ws.Range("K129").Value = 766.0908899999999
$ws.Range("L129").Value = 3276048.3
$ws.Range("M129").Value = 4233.909110000001
$ws.Range("N129").Value = -3286048.3
$ws.Range("H132").Value = 2452.138
$ws.Range("I132").Value = 2504.3462
$ws.Range("J132").Value = 1999.6666
$ws.Range("K132").Value = 7513.0386
$ws.Range("L132").Value = 5998.9998
$ws.Range("M132").Value = -4983.0386
$ws.Range("N132").Value = -11058.9998
$ws.Range("H138").Value = 2433.77
$ws.Range("I138").Value = 906.129
$ws.Range("J138").Value = 3120.1016
$ws.Range("K138").Value = 2718.387
$ws.Range("L138").Value = 9360.3048
$ws.Range("M138").Value = 2421.613
$ws.Range("N138").Value = -19640.3048
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 13261.691
$ws.Range("I32").Value = 8590.547
$ws.Range("K32").Value = 8590.547
$ws.Range("M32").Value = -8303.547
$ws.Range("H41").Value = 16215.5
$ws.Range("I41").Value = 1400
$ws.Range("J41").Value = 31031
$ws.Range("K41").Value = 1400
$ws.Range("L41").Value = 31031
$ws.Range("M41").Value = -986
$ws.Range("N41").Value = -31859
$ws.Range("H61").Value = 2554
$ws.Range("I61").Value = 2468.9412
$ws.Range("J61").Value = 4000
$ws.Range("K61").Value = 2468.9412
$ws.Range("L61").Value = 4000
$ws.Range("M61").Value = -2256.9412
$ws.Range("N61").Value = -4424
$ws.Range("H63").Value = 1996.7059
$ws.Range("I63").Value = 1996.5
$ws.Range("K63").Value = 1996.5
$ws.Range("M63").Value = -1310.5
$ws.Range("H66").Value = 1996.7059
$ws.Range("I66").Value = 1996.5
$ws.Range("K66").Value = 9982.5
$ws.Range("M66").Value = -6550.5
$ws.Range("H74").Value = 1491.1842
$ws.Range("I74").Value = 1299.5454
$ws.Range("K74").Value = 1299.5454
$ws.Range("M74").Value = -425.5454
$ws.Range("H77").Value = 1491.1842
$ws.Range("I77").Value = 1299.5454
$ws.Range("K77").Value = 6497.727
$ws.Range("M77").Value = -2129.727
$ws.Range("H102").Value = 2763.3333
$ws.Range("I102").Value = 2822.5
$ws.Range("J102").Value = 2290
$ws.Range("K102").Value = 2822.5
$ws.Range("L102").Value = 2290
$ws.Range("M102").Value = -1200.5
$ws.Range("N102").Value = -5534
$ws.Range("H122").Value = 1532.2858
$ws.Range("I122").Value = 1370.6666
$ws.Range("J122").Value = 1653.5
$ws.Range("K122").Value = 4111.9998
$ws.Range("L122").Value = 4960.5
$ws.Range("M122").Value = -1661.9998
$ws.Range("N122").Value = -9860.5
$ws.Range("H136").Value = 2554
$ws.Range("I136").Value = 2468.9412
$ws.Range("J136").Value = 4000
$ws.Range("K136").Value = 7406.823600000001
$ws.Range("L136").Value = 12000
$ws.Range("M136").Value = -4856.823600000001
$ws.Range("N136").Value = -17100
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 41924
$ws.Range("I20").Value = 121638.2
$ws.Range("J20").Value = 2066.9
$ws.Range("K20").Value = 121638.2
$ws.Range("L20").Value = 2066.9
$ws.Range("M20").Value = -121391.2
$ws.Range("N20").Value = -2560.9
$ws.Range("H86").Value = 142862200
$ws.Range("I86").Value = 200003500
$ws.Range("J86").Value = 9000
$ws.Range("K86").Value = 200003500
$ws.Range("L86").Value = 9000
$ws.Range("M86").Value = -200002377
$ws.Range("N86").Value = -11246
$ws.Range("H89").Value = 142862200
$ws.Range("I89").Value = 200003500
$ws.Range("J89").Value = 9000
$ws.Range("K89").Value = 1000017500
$ws.Range("L89").Value = 45000
$ws.Range("M89").Value = -1000011884
$ws.Range("N89").Value = -56232
$ws.Range("H105").Value = 6387
$ws.Range("I105").Value = 5105.933
$ws.Range("J105").Value = 15995
$ws.Range("K105").Value = 5105.933
$ws.Range("L105").Value = 15995
$ws.Range("M105").Value = -3358.933
$ws.Range("N105").Value = -19489
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 1996.9231
$ws.Range("J4").Value = 1996.9231
$ws.Range("L4").Value = 1996.9231
$ws.Range("N4").Value = -2220.9231
$ws.Range("H31").Value = 2065.4736
$ws.Range("I31").Value = 1624.0714
$ws.Range("J31").Value = 3301.4
$ws.Range("K31").Value = 1624.0714
$ws.Range("L31").Value = 3301.4
$ws.Range("M31").Value = -1329.0714
$ws.Range("N31").Value = -3891.4
$ws.Range("H34").Value = 2065.4736
$ws.Range("I34").Value = 1624.0714
$ws.Range("J34").Value = 3301.4
$ws.Range("K34").Value = 1624.0714
$ws.Range("L34").Value = 3301.4
$ws.Range("M34").Value = -1422.0714
$ws.Range("N34").Value = -3705.4
$ws.Range("H94").Value = 7068.8887
$ws.Range("I94").Value = 3925.5
$ws.Range("J94").Value = 9583.6
$ws.Range("K94").Value = 3925.5
$ws.Range("L94").Value = 9583.6
$ws.Range("M94").Value = -3474.5
$ws.Range("N94").Value = -10485.6
$ws.Range("H140").Value = 70342.22
$ws.Range("J140").Value = 70342.22
$ws.Range("L140").Value = 70342.22
$ws.Range("N140").Value = -80702.22
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 17277048
$ws.Range("J131").Value = 1205.7222
$ws.Range("L131").Value = 3617.1666
$ws.Range("N131").Value = -13697.1666
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H26").Value = 8625
$ws.Range("I26").Value = 8250
$ws.Range("K26").Value = 8250
$ws.Range("M26").Value = -7970
$ws.Range("H50").Value = 8625
$ws.Range("I50").Value = 8250
$ws.Range("K50").Value = 8250
$ws.Range("M50").Value = -7752
$ws.Range("H70").Value = 6510
$ws.Range("I70").Value = 6765
$ws.Range("K70").Value = 6765
$ws.Range("M70").Value = -6495
$ws.Range("H73").Value = 6510
$ws.Range("I73").Value = 6765
$ws.Range("K73").Value = 6765
$ws.Range("M73").Value = -5829
$ws.Range("H102").Value = 1690.0555
$ws.Range("I102").Value = 1509.3077
$ws.Range("J102").Value = 2160
$ws.Range("K102").Value = 1509.3077
$ws.Range("L102").Value = 2160
$ws.Range("M102").Value = 112.6922999999999
$ws.Range("N102").Value = -5404
$ws.Range("H137").Value = 35776
$ws.Range("J137").Value = 35776
$ws.Range("L137").Value = 35776
$ws.Range("N137").Value = -45976
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H2").Value = 0
$ws.Range("J2").Value = 0
$ws.Range("L2").Value = 0
$ws.Range("N2").ClearContents()
$ws.Range("H7").Value = 1921.25
$ws.Range("I7").Value = 1887.6923
$ws.Range("J7").Value = 2066.6667
$ws.Range("K7").Value = 1887.6923
$ws.Range("L7").Value = 2066.6667
$ws.Range("M7").Value = -1775.6923
$ws.Range("N7").Value = -2290.6667
$ws.Range("H40").Value = 532705.9399999999
$ws.Range("I40").Value = 777862.1
$ws.Range("K40").Value = 777862.1
$ws.Range("M40").Value = -777726.1
$ws.Range("H50").Value = 6166.6665
$ws.Range("J50").Value = 6166.6665
$ws.Range("L50").Value = 6166.6665
$ws.Range("N50").Value = -7440.6665
$ws.Range("H99").Value = 25000
$ws.Range("J99").Value = 25000
$ws.Range("L99").Value = 25000
$ws.Range("N99").Value = -30990
$ws.Range("H100").Value = 16033632
$ws.Range("I100").Value = 18705504
$ws.Range("K100").Value = 18705504
$ws.Range("M100").Value = -18704963
$ws.Range("H122").Value = 1849.375
$ws.Range("I122").Value = 1780.909
$ws.Range("K122").Value = 5342.727000000001
$ws.Range("M122").Value = -2892.727000000001
$ws.Range("H126").Value = 1921.25
$ws.Range("I126").Value = 1887.6923
$ws.Range("J126").Value = 2066.6667
$ws.Range("K126").Value = 5663.0769
$ws.Range("L126").Value = 6200.000100000001
$ws.Range("M126").Value = -3193.0769
$ws.Range("N126").Value = -11140.0001
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H102").Value = 22000
$ws.Range("J102").Value = 22000
$ws.Range("L102").Value = 22000
$ws.Range("N102").Value = -28490
